$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows continuing the time series through 2021-12-08 (row, date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newData = @"
386 44460 0 2 218.3406113537118
387 44461 0 2 218.3406113537118
388 44462 1 3 327.5109170305677
389 44463 0 3 327.5109170305677
390 44464 0 3 327.5109170305677
391 44465 0 1 109.1703056768559
392 44466 0 1 109.1703056768559
393 44467 0 1 109.1703056768559
394 44468 0 1 109.1703056768559
395 44469 0 0 0
396 44470 0 0 0
397 44471 0 0 0
398 44472 0 0 0
399 44473 0 0 0
400 44474 0 0 0
401 44475 0 0 0
402 44476 0 0 0
403 44477 0 0 0
404 44478 0 0 0
405 44479 0 0 0
406 44480 0 0 0
407 44481 0 0 0
408 44482 0 0 0
409 44483 0 0 0
410 44484 0 0 0
411 44485 0 0 0
412 44486 0 0 0
413 44487 0 0 0
414 44488 0 0 0
415 44489 0 0 0
416 44490 0 0 0
417 44491 0 0 0
418 44492 0 0 0
419 44493 0 0 0
420 44494 0 0 0
421 44495 0 0 0
422 44496 0 0 0
423 44497 0 0 0
424 44498 0 0 0
425 44499 0 0 0
426 44500 0 0 0
427 44501 0 0 0
428 44502 0 0 0
429 44503 0 0 0
430 44504 0 0 0
431 44505 0 0 0
432 44506 0 0 0
433 44507 0 0 0
434 44508 0 0 0
435 44509 0 0 0
436 44510 0 0 0
437 44511 0 0 0
438 44512 0 0 0
439 44513 0 0 0
440 44514 0 0 0
441 44515 0 0 0
442 44516 1 1 109.1703056768559
443 44517 0 1 109.1703056768559
444 44518 0 1 109.1703056768559
445 44519 0 1 109.1703056768559
446 44520 0 1 109.1703056768559
447 44521 0 1 109.1703056768559
448 44522 0 1 109.1703056768559
449 44523 0 0 0
450 44524 0 0 0
451 44525 1 1 109.1703056768559
452 44526 0 1 109.1703056768559
453 44527 0 1 109.1703056768559
454 44528 0 1 109.1703056768559
455 44529 0 1 109.1703056768559
456 44530 1 2 218.3406113537118
457 44531 0 2 218.3406113537118
458 44532 0 1 109.1703056768559
459 44533 0 1 109.1703056768559
460 44534 0 1 109.1703056768559
461 44535 0 1 109.1703056768559
462 44536 0 1 109.1703056768559
463 44537 0 0 0
464 44538 0 0 0
"@

# Copy column-A date formatting/style from the last existing data row (385) down across
# the whole new block so the appended dates render the same as the pre-existing ones.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)

foreach ($line in ($newData -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\s+"
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = [double]$parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [double]$parts[4]
}
